# Sample Project / Main.xlsx - "Rules" sheet, row 11 ("R40" rule).
# B11 previously held the text "R40"; it is retyped to the text "1".
# A leading apostrophe forces Excel to keep the numeric-looking input as
# text (matching the workbook's original shared-string "t=s" cell type)
# instead of silently coercing it to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "'1"
